$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to Text format before assigning values so that
# values such as "65.730.53" (dotted thousands) or "1.00" / "0.0260" are not
# auto-coerced into numbers (which would corrupt or reformat them), matching
# the source data which stores these as plain text strings.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.730.53'
$ws.Range("E2").Value = '  -0.85%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.672.38'
$ws.Range("E3").Value = '  -0.43%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.32'
$ws.Range("E5").Value = '  -1.34%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.79'
$ws.Range("E6").Value = '  -1.53%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.622'
$ws.Range("E8").Value = '  +4.55%  '

$ws.Range("E9").Value = '  +1.45%  '

$ws.Range("E10").Value = '  -1.38%  '

$ws.Range("E11").Value = '  -2.85%  '

$ws.Range("E12").Value = '  -0.37%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '29.54'
$ws.Range("E13").Value = '  -3.67%  '

$ws.Range("E14").Value = '  -7.10%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.152.47'
$ws.Range("E15").Value = '  -0.53%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.593.50'
$ws.Range("E16").Value = '  -0.79%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.647.22'
$ws.Range("E17").Value = '  -1.38%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.84'
$ws.Range("E18").Value = '  +0.69%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.82'
$ws.Range("E19").Value = '  -1.84%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.64'
$ws.Range("E20").Value = '  +1.90%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '351.70'
$ws.Range("E21").Value = '  -3.07%  '

$ws.Range("E22").Value = '  +0.01%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.95'
$ws.Range("E23").Value = '  -0.48%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000111'
$ws.Range("E24").Value = '  +2.38%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.82'
$ws.Range("E25").Value = '  +0.66%  '

$ws.Range("E26").Value = '  -2.72%  '

$ws.Range("E27").Value = '  -3.92%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.58'
$ws.Range("E28").Value = '  -5.33%  '

$ws.Range("E29").Value = '  -1.61%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.05%  '

$ws.Range("E31").Value = '  -2.57%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '526.88'
$ws.Range("E32").Value = '  -2.96%  '

$ws.Range("E33").Value = '  -2.20%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.55'
$ws.Range("E34").Value = '  -0.99%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.50'
$ws.Range("E35").Value = '  -1.49%  '

$ws.Range("E36").Value = '  -2.22%  '

$ws.Range("E37").Value = '  -1.82%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '159.84'
$ws.Range("E38").Value = '  -1.99%  '

$ws.Range("E39").Value = '  -0.05%  '

$ws.Range("E40").Value = '  -4.37%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '42.73'
$ws.Range("E42").Value = '  +0.67%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '165.19'
$ws.Range("E43").Value = '  -3.34%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.11'
$ws.Range("E44").Value = '  -4.01%  '

$ws.Range("E45").Value = '  -1.77%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0612'
$ws.Range("E46").Value = '  -1.37%  '

$ws.Range("E47").Value = '  -0.39%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0260'
$ws.Range("E48").Value = '  -2.18%  '

$ws.Range("E49").Value = '  -3.13%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.102'
$ws.Range("E50").Value = '  +2.38%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.23'
$ws.Range("E51").Value = '  +0.75%  '
